# Update numeric values in Sheet1 cells as produced by a re-run of the
# RandomForest imputation algorithm (commit message: "Update Name of Algo").
# The diff only touches the <v> value of specific existing numeric cells in
# columns B, C, and D; no structural/style changes are required.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B3" = 5.500900000000002
    "D6" = -8.300200000000004
    "B14" = 6.9375
    "C15" = -13.37269999999999
    "B16" = 5.4655
    "D18" = -8.702199999999999
    "D19" = -8.14589999999999
    "B21" = 8.538199999999996
    "C21" = -11.76810000000001
    "C22" = -11.8136
    "B23" = 9.0984
    "C24" = -13.56159999999999
    "B25" = 6.704199999999997
    "B26" = 5.540100000000003
    "C27" = -12.83539999999999
    "C28" = -13.5122
    "B29" = 5.788399999999998
    "D35" = -8.554099999999988
    "C36" = -13.41880000000001
    "C39" = -13.40490000000001
    "B40" = 9.030899999999992
    "D44" = -7.404299999999996
    "C45" = -13.70759999999999
    "D47" = -7.085099999999998
    "C48" = -11.8485
    "C49" = -14.11129999999999
    "D50" = -8.625400000000003
    "D51" = -8.619200000000003
    "C52" = -10.9747
    "D52" = -7.015399999999996
    "B53" = 6.803999999999997
    "C53" = -11.1333
    "C54" = -13.52120000000001
    "D55" = -8.857000000000001
    "B57" = 5.305100000000001
    "C57" = -13.59599999999999
    "D57" = -9.006599999999993
    "D58" = -8.265100000000011
    "B59" = 4.9254
    "D64" = -7.859099999999993
    "B65" = 6.990100000000004
    "D66" = -7.735299999999992
    "B69" = 5.111899999999998
    "C70" = -11.5226
    "C71" = -10.72040000000001
    "B79" = 9.509100000000004
    "D80" = -7.720099999999993
    "B83" = 5.454799999999999
    "D83" = -9.158099999999994
    "C86" = -14.12339999999999
    "C87" = -12.98709999999999
    "C89" = -12.9682
    "B91" = 5.505700000000003
    "D92" = -7.112299999999997
    "B93" = 5.3236
    "D94" = -6.113900000000001
    "D96" = -8.431400000000007
    "D97" = -8.466699999999989
    "B100" = 5.012000000000004
    "C101" = -11.3319
    "D101" = -7.096300000000005
    "B103" = 5.290100000000006
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
